# Update the "Förändrad" date column (C) for rows 2-7 from 2023-09-06 to 2023-09-14
# (serial date 45175 -> 45183), matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45183
}
